# Objetivo_3.xlsx — update current-sensor units from mA to A.
#
# The label for the current measurement changes from "Corrente (mA)" to
# "Corrente (A)" (B7), the measured current value changes from 5 (mA) to
# 0.52 (A) (B8), and the two formulas that used to convert mA -> A by
# dividing by 1000 (J7, J19) are simplified since B8 is now already in
# amps. The active selection also moves from F18 to J8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current label: "Corrente (mA)" -> "Corrente (A)"
$ws.Range("B7").Value = "Corrente (A)"

# Current value: 5 mA -> 0.52 A
$ws.Range("B8").Value = 0.52

# Formulas no longer need to divide B8 by 1000 (B8 is already in amps)
$ws.Range("J7").Formula = "=(B11-(J15*(B8)))/F6"
$ws.Range("J19").Formula = "=J7*B8"

# Update the active cell selection on the sheet
[void]$ws.Range("J8").Select()
